$wb = $excel.ActiveWorkbook

# Rename sheets (task order identifiers updated to new timestamps)
$wb.Worksheets.Item(1).Name = "GNG_TO-16504778218458996"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778244478853"
$wb.Worksheets.Item(3).Name = "RS_TO-1650477824448882"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778244959104"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778245598798"

# Sheet 1 (GNG) - update stimulus file timestamps
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650477821801879.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778218278792.csv"
$ws1.Range("B4").Value = "go_stims-16504778218288815.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778218438785.csv"

# Sheet 2 (NB) - update stimulus file timestamps
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16504778225988781.csv"
$ws2.Range("B3").Value = "OB-16504778227078753.csv"
$ws2.Range("B4").Value = "TB-16504778236099095.csv"
$ws2.Range("B5").Value = "OB-1650477822992913.csv"
$ws2.Range("B6").Value = "ZB-match_1-16504778220438795.csv"
$ws2.Range("B7").Value = "ZB-match_5-16504778222678788.csv"
$ws2.Range("B8").Value = "TB-16504778233718762.csv"
$ws2.Range("B9").Value = "ZB-match_0-1650477822219878.csv"
$ws2.Range("B10").Value = "TB-1650477824428908.csv"

# Sheet 4 (TOL) - update stimulus file timestamps
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778244638805.csv"
$ws4.Range("B3").Value = "ZM_stims-1650477824450883.csv"
$ws4.Range("B4").Value = "MM_stims-16504778244799113.csv"
$ws4.Range("B5").Value = "ZM_stims-1650477824464878.csv"
$ws4.Range("B6").Value = "MM_stims-16504778244959104.csv"
$ws4.Range("B7").Value = "ZM_stims-1650477824480878.csv"

# Sheet 5 (vSAT) - update stimulus file timestamps
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16504778245128775.csv"
$ws5.Range("B3").Value = "SAT_stims-165047782449888.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778245279112.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778245439112.csv"
